$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.211.86'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.002.94'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.34'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("E6").Value = '  +2.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.31'
$ws.Range("E7").Value = '  +3.33%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("E10").Value = '  +2.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.08'
$ws.Range("E12").Value = '  +7.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.83'
$ws.Range("E13").Value = '  +7.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.850'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.297.15'
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.47'
$ws.Range("E16").Value = '  +3.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.016.67'
$ws.Range("E17").Value = '  +3.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.188.27'
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.46'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0866'
$ws.Range("E20").Value = '  +2.26%  '
$ws.Range("E21").Value = '  +3.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.92'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("E26").Value = '  +3.41%  '
$ws.Range("E27").Value = '  +2.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.11'
$ws.Range("E28").Value = '  +2.34%  '
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.34'
$ws.Range("E30").Value = '  +13.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.121'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0659'
$ws.Range("E33").Value = '  +7.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.53'
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  +5.62%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +2.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.44'
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  +2.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.19'
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.82'
$ws.Range("E44").Value = '  +7.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.29'
$ws.Range("E45").Value = '  +3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.383.97'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  +3.75%  '
$ws.Range("E49").Value = '  +14.87%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.48'
$ws.Range("E51").Value = '  +5.50%  '
